# Generate Report for Handback
# - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#   (this text lives once in the shared-string table and is referenced by the
#   Overview sheet's zh-cn/de-de columns as well as by the Status column on
#   both the zh-cn and de-de detail sheets, so a single Range write per
#   sheet/column is enough to flip every occurrence).
# - The de-de handback round-trip finished: a fresh "Latest Handback DateTime"
#   is stamped and the stale "handback file is not latest" error is cleared.
# - The zh-cn handback timestamp is likewise refreshed.
# - Column widths for the Status / Error Detail columns are refreshed to fit
#   the new text.

$wb = $excel.ActiveWorkbook
$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Status text (Overview!E2:F2, zh-cn!C2, de-de!C2) ---
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsZhCn.Range("C2").Value = $newStatus
$wsDeDe.Range("C2").Value = $newStatus

# --- zh-cn: refresh handback datetime, clear stale error ---
$wsZhCn.Range("K2").Value = "2016-08-26 02:49:43"
$wsZhCn.Range("P2").Value = ""

# --- de-de: refresh handback datetime, clear stale error ---
$wsDeDe.Range("K2").Value = "2016-08-26 02:49:50"
$wsDeDe.Range("P2").Value = ""

# --- Column width refresh to fit the new Status / cleared Error Detail text ---
$wsOverview.Range("E1").ColumnWidth = 29.1
$wsOverview.Range("F1").ColumnWidth = 29.1
$wsZhCn.Range("C1").ColumnWidth = 29.1
$wsZhCn.Range("P1").ColumnWidth = 12.8
$wsDeDe.Range("C1").ColumnWidth = 29.1
$wsDeDe.Range("P1").ColumnWidth = 12.8
